$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.776.93'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.41%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.980.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.18%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '540.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.93%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.17'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.66%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.568'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.44%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.992.44'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.20%  '

$ws.Range('E10').Value = '  -4.05%  '

$ws.Range('E11').Value = '  -7.33%  '

$ws.Range('E12').Value = '  -4.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.499.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.19%  '

$ws.Range('E14').Value = '  -2.36%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '61.792.43'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.36%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.22%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.982.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.34%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000147'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.80%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.37%  '

$ws.Range('E20').Value = '  -3.73%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.14%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.27%  '

$ws.Range('E23').Value = '  +0.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.23%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.472'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.75%  '

$ws.Range('E26').Value = '  -5.34%  '

$ws.Range('E27').Value = '  -2.76%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.53%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0942'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.04%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.12'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.94%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.73%  '

$ws.Range('E33').Value = '  -5.48%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '159.66'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.42%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.93'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.99%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.57'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.28%  '

$ws.Range('E37').Value = '  -5.41%  '

$ws.Range('E38').Value = '  -6.84%  '

$ws.Range('E39').Value = '  -8.89%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.61%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.423.42'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.27%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.90'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.99%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.63%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.672'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.78%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0590'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.86%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.76%  '

$ws.Range('E47').Value = '  +0.14%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0245'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.05%  '

$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.82'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.63%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0953'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.35%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '267.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.29%  '
